$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing data row (row 2) in place: date + start/end hours change.
$ws.Range("A2").Value = 45063
$ws.Range("B2").Value = 15
$ws.Range("C2").Value = 17

# Insert four new rows below it (rows 3-6), inheriting row 2's formatting
# (date format on column A, general/applied format on columns B and C).
$ws.Rows("3:6").Insert()

# Row 3: the day after the (old) first row's original date.
$ws.Range("A3").Value = 45064
$ws.Range("B3").Value = 16
$ws.Range("C3").Value = 18

# Row 4: the original date/hours that used to live in row 2.
$ws.Range("A4").Value = 45402
$ws.Range("B4").Value = 16
$ws.Range("C4").Value = 18

# Rows 5 and 6 stay empty (only A/C keep formatting, no B cell at all).
$ws.Range("B5").Clear()
$ws.Range("B6").Clear()

# Match the saved selection state from the edit.
$ws.Range("F15").Select()
